$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# Add the new row of mail-log data
$ws.Range("A7").Value = "Vragen over samenwerking"
$ws.Range("B7").Value = "mailmind.test@zohomail.eu"
$ws.Range("C7").Value = "Kunnen we samenwerken aan een nieuw project?"
$ws.Range("D7").Value = "Overig"
$ws.Range("F7").Value = "2025-06-19 10:58:11"
$ws.Range("G7").Value = "Nee"

# Update the Dashboard summary count for "Overig"
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 4

# Extend the conditional formatting ranges to cover the new row
$catRules = $ws.Range("D2:D6").FormatConditions
$catRules.Item(1).ModifyAppliesToRange($ws.Range("D2:D7"))

$answeredRules = $ws.Range("G2:G6").FormatConditions
$answeredRules.Item(1).ModifyAppliesToRange($ws.Range("G2:G7"))
